# Daily attendance processing - 2025-11-17 20:24:39
# Normalizes the "Recorded By" (column G) values: for rows where the
# recorder list does not already start with dnasr281@gmail.com, swap the
# first two comma-separated entries so the human reviewer's address is
# listed ahead of the automated "System" entry.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 7)
    $val = $cell.Value2

    if ($val -ne $null) {
        $parts = $val -split ", "

        if ($parts.Length -ge 2) {
            if ($parts[0] -ne "dnasr281@gmail.com") {
                $first = $parts[0]
                $second = $parts[1]
                $parts[0] = $second
                $parts[1] = $first
                $cell.Value2 = $parts -join ", "
            }
        }
    }
}
